{"js": "const styles = context.document.getStyles();\n\n// Remove the custom \"Abstract Title\" paragraph style entirely.\nconst abstractTitleStyle = styles.getByNameOrNullObject(\"Abstract Title\");\nabstractTitleStyle.load(\"isNullObject\");\nawait context.sync();\nif (!abstractTitleStyle.isNullObject) {\n  abstractTitleStyle.delete();\n  await context.sync();\n}\n\n// Update the \"Abstract\" style: space-before goes from 100 twips (5pt) to 300 twips (15pt).\nconst abstractStyle = styles.getByNameOrNullObject(\"Abstract\");\nabstractStyle.load(\"isNullObject\");\nawait context.sync();\nif (!abstractStyle.isNullObject) {\n  abstractStyle.paragraphFormat.spaceBefore = 15;\n  await context.sync();\n}\n\n// Remove the custom \"Footnote Block Text\" paragraph style entirely.\nconst footnoteBlockTextStyle = styles.getByNameOrNullObject(\"Footnote Block Text\");\nfootnoteBlockTextStyle.load(\"isNullObject\");\nawait context.sync();\nif (!footnoteBlockTextStyle.isNullObject) {\n  footnoteBlockTextStyle.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the custom \"Abstract Title\" paragraph style entirely.\ntry {\n    $d.Styles(\"Abstract Title\").Delete()\n} catch {}\n\n# Update the \"Abstract\" style: space-before goes from 100 twips (5pt) to 300 twips (15pt).\ntry {\n    $d.Styles(\"Abstract\").ParagraphFormat.SpaceBefore = 15\n} catch {}\n\n# Remove the custom \"Footnote Block Text\" paragraph style entirely.\ntry {\n    $d.Styles(\"Footnote Block Text\").Delete()\n} catch {}\n"}
